$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet had two weekly entries (Primera/Segunda) at rows 96-97 (week of
# 44225) and rows 98-99 (week of 44236). A new week (44890) was recorded and
# inserted ahead of the 44236 week, so:
#   - rows 98-99 (44236 week) shift down to rows 100-101, unchanged
#   - rows 98-99 are freed up and now hold what used to be the 44225 week
#     (previously at rows 96-97)
#   - rows 96-97 are overwritten with the new 44890 week values
# ---------------------------------------------------------------------------

# Insert two blank rows at row 98, pushing the old rows 98-99 down to 100-101
$ws.Rows.Item(98).Insert()
$ws.Rows.Item(98).Insert()

# New values for row 96 (was: D=44225, M=100, N=3000, O=3500, P=3250, S=1625)
$ws.Range("D96").Value = 44890
$ws.Range("M96").Value = 200
$ws.Range("N96").Value = 4000
$ws.Range("O96").Value = 4500
$ws.Range("P96").Value = 4250
$ws.Range("S96").Value = 2125

# New values for row 97 (was: D=44225, M=50, N=2500, O=2500, P=2500, S=1250)
$ws.Range("D97").Value = 44890
$ws.Range("M97").Value = 100
$ws.Range("N97").Value = 3500
$ws.Range("O97").Value = 3500
$ws.Range("P97").Value = 3500
$ws.Range("S97").Value = 1750

# Row 98 now holds what used to be in row 96 (44225 week, Primera)
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = "Vega Monumental Concepción"
$ws.Range("C98").Value = "Bíobío"
$ws.Range("D98").Value = 44225
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100101
$ws.Range("H98").Value = "Berries"
$ws.Range("I98").Value = 100101001
$ws.Range("J98").Value = "Arándano (blue)"
$ws.Range("K98").Value = "Sin especificar"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 100
$ws.Range("N98").Value = 3000
$ws.Range("O98").Value = 3500
$ws.Range("P98").Value = 3250
$ws.Range("Q98").Value = "$/bandeja 2 kilos"
$ws.Range("R98").Value = "Región de Ñuble"
$ws.Range("S98").Value = 1625
$ws.Range("T98").Value = 2

# Row 99 now holds what used to be in row 97 (44225 week, Segunda)
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = "Vega Monumental Concepción"
$ws.Range("C99").Value = "Bíobío"
$ws.Range("D99").Value = 44225
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100101
$ws.Range("H99").Value = "Berries"
$ws.Range("I99").Value = 100101001
$ws.Range("J99").Value = "Arándano (blue)"
$ws.Range("K99").Value = "Sin especificar"
$ws.Range("L99").Value = "Segunda"
$ws.Range("M99").Value = 50
$ws.Range("N99").Value = 2500
$ws.Range("O99").Value = 2500
$ws.Range("P99").Value = 2500
$ws.Range("Q99").Value = "$/bandeja 2 kilos"
$ws.Range("R99").Value = "Región de Ñuble"
$ws.Range("S99").Value = 1250
$ws.Range("T99").Value = 2

# Rows 100-101 already hold the shifted-down former rows 98-99 (44236 week)
# with correct values/styles thanks to the Insert() above; nothing else to do.
